# Updates cryptos list values to match the latest scrape (GitHub Actions run).
# Most D-column prices look like plain numbers ("309.24") but the sheet stores
# them as text (inline strings), so we force NumberFormat="@" before writing and
# reset the style back to Normal afterwards to avoid leaving a text format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.125.03'
$ws.Range("E2").Value = '  +2.56%  '
$ws.Range("D3").Value = '2.348.82'
$ws.Range("E3").Value = '  +7.06%  '
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("E7").Value = '  +3.66%  '
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.05'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.64%  '
$ws.Range("E14").Value = '  +2.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +11.44%  '
$ws.Range("D16").Value = '2.705.19'
$ws.Range("E16").Value = '  +6.99%  '
$ws.Range("D17").Value = '2.489.05'
$ws.Range("E17").Value = '  +12.50%  '
$ws.Range("D18").Value = '43.062.43'
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("E19").Value = '  +4.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +13.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '252.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.88%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.35%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0920'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.93'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.07%  '
$ws.Range("E36").Value = '  +3.58%  '
$ws.Range("E37").Value = '  +6.90%  '
$ws.Range("E38").Value = '  +6.49%  '
$ws.Range("E39").Value = '  -1.65%  '
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.89%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.86%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +17.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.231'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.89%  '
$ws.Range("E50").Value = '  +1.17%  '
$ws.Range("E51").Value = '  +3.51%  '
